# Update "想去人数" (want-to-go count) figures across sheets, reflecting a
# newer data pull (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 287
$ws1.Range("F4").Value = 230
$ws1.Range("F5").Value = 1725
$ws1.Range("F7").Value = 630
$ws1.Range("F8").Value = 136
$ws1.Range("F9").Value = 602
$ws1.Range("F10").Value = 61
$ws1.Range("F13").Value = 163

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1947
$ws3.Range("F5").Value = 83

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1947
$ws4.Range("F5").Value = 287
$ws4.Range("F6").Value = 83
$ws4.Range("F12").Value = 230
$ws4.Range("F16").Value = 1725
$ws4.Range("F21").Value = 630
$ws4.Range("F23").Value = 136
$ws4.Range("F24").Value = 602
$ws4.Range("F25").Value = 61
$ws4.Range("F31").Value = 163

$wb.Save()
